$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "iPhone 14 Pro"
$ws.Range("B7").Value = 20000.0
$ws.Range("C7").Value = 20.0
$ws.Range("D7").Value = "Electronics"
